$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 106: the R refresh re-pulled this bar with a corrected timestamp
#     (midnight snapshot instead of an intraday one) and revised OHLC values.
#     close (F), adj_close (G) and ticker (H) were already correct. ---
$ws.Range("A106").Value = 45467.2916666667
$ws.Range("C106").Value = 1.89999997615814
$ws.Range("D106").Value = 1.89999997615814
$ws.Range("E106").Value = 1.89999997615814

# --- Row 107: brand-new bar appended by the R script. ---

# A107 needs the same date/time number format as the rest of column A, so
# copy A106's formatting down before writing the new timestamp.
$ws.Range("A106").Copy()
$ws.Range("A107").PasteSpecial(-4122)
$ws.Range("A107").Value = 45468.2925694444

$ws.Range("B107").Value = 900
$ws.Range("C107").Value = 1.91999995708466
$ws.Range("D107").Value = 1.91999995708466
$ws.Range("E107").Value = 1.91999995708466
$ws.Range("F107").Value = 1.91999995708466

# adj_close (column G) is stored as literal text in this sheet, not as a
# number, so the numeric-looking value has to be forced into a text cell.
# Stage it in a scratch cell (with NumberFormat "@" so Excel doesn't
# auto-convert it back to a number), copy just that cell's value+format
# into G107, then remove the scratch cell/column again so no trace of it
# is left behind in the saved sheet.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "1.91999995708466"
$ws.Range("Z1").Copy()
$ws.Range("G107").PasteSpecial(-4163)
$ws.Range("Z1").EntireColumn.Delete()

$ws.Range("H107").Value = "KK.MI"
